$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 567.35297
$ws.Range("I28").Value = 587.1875
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 587.1875
$ws.Range("L28").Value = 250
$ws.Range("M28").Value = -102.1875
$ws.Range("N28").Value = -1220
# Row 92
$ws.Range("H92").Value = 293.08334
$ws.Range("I92").Value = 192.45454
$ws.Range("J92").Value = 1400
$ws.Range("K92").Value = 192.45454
$ws.Range("L92").Value = 1400
$ws.Range("M92").Value = 1055.54546
$ws.Range("N92").Value = -3896
# Row 99
$ws.Range("H99").Value = 1028.4546
$ws.Range("I99").Value = 807.4286
$ws.Range("J99").Value = 1415.25
$ws.Range("K99").Value = 2422.2858
$ws.Range("L99").Value = 4245.75
$ws.Range("M99").Value = -924.2857999999997
$ws.Range("N99").Value = -7241.75
# Row 101
$ws.Range("H101").Value = 3350.4285
$ws.Range("I101").Value = 5217
$ws.Range("J101").Value = 861.6667
$ws.Range("K101").Value = 15651
$ws.Range("L101").Value = 2585.0001
$ws.Range("M101").Value = -14029
$ws.Range("N101").Value = -5829.0001
# Row 113
$ws.Range("H113").Value = 2647
$ws.Range("I113").Value = 1951.375
$ws.Range("J113").Value = 3760
$ws.Range("K113").Value = 1951.375
$ws.Range("L113").Value = 3760
$ws.Range("M113").Value = 1302.625
$ws.Range("N113").Value = -10268
# Row 115
$ws.Range("H115").Value = 4980
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 4980
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = ""
$ws.Range("M115").Value = 14940
$ws.Range("N115").Value = -18074
# Row 125
$ws.Range("H125").Value = 112001.445
$ws.Range("I125").Value = 200572.6
$ws.Range("J125").Value = 1287.5
$ws.Range("K125").Value = 1805153.4
$ws.Range("L125").Value = 11587.5
$ws.Range("M125").Value = -1802693.4
$ws.Range("N125").Value = -16507.5
# Row 129
$ws.Range("H129").Value = 1210.0588
$ws.Range("I129").Value = 626.3333
$ws.Range("J129").Value = 1528.4546
$ws.Range("K129").Value = 1878.9999
$ws.Range("L129").Value = 4585.3638
$ws.Range("M129").Value = 3121.0001
$ws.Range("N129").Value = -14585.3638
# Row 138
$ws.Range("H138").Value = 3575844
$ws.Range("I138").Value = 8334962
$ws.Range("J138").Value = 6505.375
$ws.Range("K138").Value = 25004886
$ws.Range("L138").Value = 19516.125
$ws.Range("M138").Value = -24999746
$ws.Range("N138").Value = -29796.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = ""
$ws.Range("N24").Value = 0
# Row 25
$ws.Range("H25").Value = 25500
$ws.Range("I25").Value = 25500
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 25500
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = -25098
# Row 32
$ws.Range("H32").Value = 17779.682
$ws.Range("I32").Value = 17436.777
$ws.Range("J32").Value = 25495
$ws.Range("K32").Value = 17436.777
$ws.Range("L32").Value = 25495
$ws.Range("M32").Value = -17149.777
$ws.Range("N32").Value = -26069
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = ""
# Row 74
$ws.Range("H74").Value = 1207.4
$ws.Range("I74").Value = 1094.8
$ws.Range("J74").Value = 1320
$ws.Range("K74").Value = 1094.8
$ws.Range("L74").Value = 1320
$ws.Range("M74").Value = -220.8
$ws.Range("N74").Value = -3068
# Row 77
$ws.Range("H77").Value = 1207.4
$ws.Range("I77").Value = 1094.8
$ws.Range("J77").Value = 1320
$ws.Range("K77").Value = 5474
$ws.Range("L77").Value = 6600
$ws.Range("M77").Value = -1106
$ws.Range("N77").Value = -15336
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = ""
$ws.Range("N100").Value = 0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 24
$ws.Range("H24").Value = 1264
$ws.Range("I24").Value = 1391.3334
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 1391.3334
$ws.Range("L24").Value = 500
$ws.Range("M24").Value = -1156.3334
$ws.Range("N24").Value = -970
# Row 29
$ws.Range("H29").Value = 866.6667
$ws.Range("I29").Value = 866.6667
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 866.6667
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -577.6667
# Row 32
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = ""
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = ""
$ws.Range("N34").Value = 0
# Row 37
$ws.Range("H37").Value = 2932.4
$ws.Range("I37").Value = 1178
$ws.Range("J37").Value = 9950
$ws.Range("K37").Value = 1178
$ws.Range("L37").Value = 9950
$ws.Range("M37").Value = -1041
$ws.Range("N37").Value = -10224
# Row 94
$ws.Range("H94").Value = 1296.6818
$ws.Range("I94").Value = 1038.5625
$ws.Range("J94").Value = 1985
$ws.Range("K94").Value = 1038.5625
$ws.Range("L94").Value = 1985
$ws.Range("M94").Value = -587.5625
$ws.Range("N94").Value = -2887

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 18
$ws.Range("H18").Value = 23250
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 23250
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 23250
$ws.Range("N18").Value = -23710
# Row 31
$ws.Range("H31").Value = 29415030
$ws.Range("I31").Value = 43480852
$ws.Range("J31").Value = 4681.909
$ws.Range("K31").Value = 43480852
$ws.Range("L31").Value = 4681.909
$ws.Range("M31").Value = -43480557
$ws.Range("N31").Value = -5271.909
# Row 34
$ws.Range("H34").Value = 29415030
$ws.Range("I34").Value = 43480852
$ws.Range("J34").Value = 4681.909
$ws.Range("K34").Value = 43480852
$ws.Range("L34").Value = 4681.909
$ws.Range("M34").Value = -43480650
$ws.Range("N34").Value = -5085.909
# Row 109
$ws.Range("H109").Value = 20000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 20000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22080

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1299.8572
$ws.Range("I68").Value = 1332.6666
$ws.Range("J68").Value = 1286.7333
$ws.Range("K68").Value = 3997.9998
$ws.Range("L68").Value = 3860.199900000001
$ws.Range("M68").Value = -3186.9998
$ws.Range("N68").Value = -5482.199900000001
# Row 71
$ws.Range("H71").Value = 1299.8572
$ws.Range("I71").Value = 1332.6666
$ws.Range("J71").Value = 1286.7333
$ws.Range("K71").Value = 11993.9994
$ws.Range("L71").Value = 11580.5997
$ws.Range("M71").Value = -7937.999400000001
$ws.Range("N71").Value = -19692.5997
# Row 107
$ws.Range("H107").Value = 694.0909
$ws.Range("I107").Value = 1467.6666
$ws.Range("J107").Value = 404
$ws.Range("K107").Value = 4402.9998
$ws.Range("L107").Value = 1212
$ws.Range("M107").Value = -2482.9998
$ws.Range("N107").Value = -5052
# Row 113
$ws.Range("H113").Value = 656.13635
$ws.Range("I113").Value = 467.77777
$ws.Range("J113").Value = 786.53845
$ws.Range("K113").Value = 1403.33331
$ws.Range("L113").Value = 2359.61535
$ws.Range("M113").Value = 766.66669
$ws.Range("N113").Value = -6699.61535
# Row 121
$ws.Range("H121").Value = 63651.625
$ws.Range("I121").Value = 1590
$ws.Range("J121").Value = 84338.836
$ws.Range("K121").Value = 4770
$ws.Range("L121").Value = 253016.508
$ws.Range("M121").Value = -3460
$ws.Range("N121").Value = -255636.508
# Row 122
$ws.Range("H122").Value = 771.7778
$ws.Range("I122").Value = 491.77777
$ws.Range("J122").Value = 1331.7778
$ws.Range("K122").Value = 4425.99993
$ws.Range("L122").Value = 11986.0002
$ws.Range("M122").Value = -1975.99993
$ws.Range("N122").Value = -16886.0002
# Row 131
$ws.Range("H131").Value = 873.21
$ws.Range("I131").Value = 485
$ws.Range("J131").Value = 889.38544
$ws.Range("K131").Value = 1455
$ws.Range("L131").Value = 2668.15632
$ws.Range("M131").Value = 3585
$ws.Range("N131").Value = -12748.15632

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6121.4707
$ws.Range("I70").Value = 5755.2144
$ws.Range("J70").Value = 6377.85
$ws.Range("K70").Value = 5755.2144
$ws.Range("L70").Value = 6377.85
$ws.Range("M70").Value = -5485.2144
$ws.Range("N70").Value = -6917.85
# Row 73
$ws.Range("H73").Value = 6121.4707
$ws.Range("I73").Value = 5755.2144
$ws.Range("J73").Value = 6377.85
$ws.Range("K73").Value = 5755.2144
$ws.Range("L73").Value = 6377.85
$ws.Range("M73").Value = -4819.2144
$ws.Range("N73").Value = -8249.85
# Row 132
$ws.Range("H132").Value = 2782.5652
$ws.Range("I132").Value = 2085.5454
$ws.Range("J132").Value = 3421.5
$ws.Range("K132").Value = 6256.6362
$ws.Range("L132").Value = 10264.5
$ws.Range("M132").Value = -3726.6362
$ws.Range("N132").Value = -15324.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = ""
$ws.Range("N3").Value = 0
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = ""
$ws.Range("N15").Value = 0
# Row 82
$ws.Range("H82").Value = 2742.1428
$ws.Range("I82").Value = 1956.5555
$ws.Range("J82").Value = 4156.2
$ws.Range("K82").Value = 1956.5555
$ws.Range("L82").Value = 4156.2
$ws.Range("M82").Value = -1595.5555
$ws.Range("N82").Value = -4878.2
# Row 85
$ws.Range("H85").Value = 2742.1428
$ws.Range("I85").Value = 1956.5555
$ws.Range("J85").Value = 4156.2
$ws.Range("K85").Value = 1956.5555
$ws.Range("L85").Value = 4156.2
$ws.Range("M85").Value = -708.5554999999999
$ws.Range("N85").Value = -6652.2
# Row 100
$ws.Range("H100").Value = 5514
$ws.Range("I100").Value = 7835.375
$ws.Range("J100").Value = 1799.8
$ws.Range("K100").Value = 7835.375
$ws.Range("L100").Value = 1799.8
$ws.Range("M100").Value = -7294.375
$ws.Range("N100").Value = -2881.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 798
$ws.Range("I96").Value = 866.6667
$ws.Range("J96").Value = 695
$ws.Range("K96").Value = 866.6667
$ws.Range("L96").Value = 695
$ws.Range("M96").Value = 506.3333
$ws.Range("N96").Value = -3441
# Row 109
$ws.Range("H109").Value = 48377
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 48377
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 48377
$ws.Range("N109").Value = -51151
